$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data (row 8), continuing the training log
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 42604.890416666669

$ws.Range("B8").Value = "Bag"

$ws.Range("C8").Value = 4034
$ws.Range("D8").Value = 4529
$ws.Range("E8").Value = 571
$ws.Range("F8").Value = 49
$ws.Range("G8").Value = 61
$ws.Range("H8").Value = 44
$ws.Range("I8").Value = 54
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 8
$ws.Range("L8").Value = 20
$ws.Range("M8").Value = 80
